$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values (translations) ---

# Row 1 header cells: M1/N1 values unchanged ("es"/"he_IL"); only style changes (below)

# Row 2
$ws.Range("D2").Value = "test"
$ws.Range("F2").Value = "test"
$ws.Range("H2").Value = "test"
$ws.Range("I2").Value = "test"
$ws.Range("K2").Value = "test"
$ws.Range("L2").Value = "test"
$ws.Range("M2").Value = "test"
$ws.Range("N2").Value = "test"

# Row 3
$ws.Range("D3").Value = "test"
$ws.Range("F3").Value = "test"
$ws.Range("H3").Value = "test"
$ws.Range("I3").Value = "test"
$ws.Range("K3").Value = "test"
$ws.Range("L3").Value = "test"
$ws.Range("M3").Value = "test"
$ws.Range("N3").Value = "test"

# Row 4
$ws.Range("D4").Value = "test"
$ws.Range("F4").Value = "test"
$ws.Range("H4").Value = "test"
$ws.Range("I4").Value = "test"
$ws.Range("K4").Value = "test"
$ws.Range("L4").Value = "test"
$ws.Range("M4").Value = "test"
$ws.Range("N4").Value = "test"

# --- Apply formatting (style index 1, same as column A) to the translated cells ---
# Use copy/paste-special of formats from a cell that already carries style "1" (A1)
# so we reuse the existing cellXfs entry instead of creating a new one.

$ws.Range("A1").Copy()
$targetRanges = @(
  "M1", "N1",
  "D2", "F2", "H2", "I2", "K2", "L2", "M2", "N2",
  "D3", "F3", "H3", "I3", "K3", "L3", "M3", "N3",
  "D4", "F4", "H4", "I4", "K4", "L4", "M4", "N4"
)
foreach ($addr in $targetRanges) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}

# --- J2:J4 lose their explicit style (revert to default / no "s" attribute) ---
$ws.Range("J2:J4").ClearFormats()

Write-Output "done"
